# Generate Report for Handback
# The localization pipeline has handed back the 7de06e11-... file (and its
# duplicate a5786877-...) in both the zh-cn and de-de target languages.
# Update the workbook to reflect the new "Handed back" status and record
# the target/handback files + timestamps for each language sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update status columns (zh-cn / de-de) for the
# 7de06e11-3939-4d03-b6a1-b5617c4c3b14.md row (row 3)
# ---------------------------------------------------------------------
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 (7de06e11-...) and row 4 (a5786877-..., content
# duplicate of row 3) now have a handback target/file recorded.
# ---------------------------------------------------------------------
$wsZhCn.Range("C3").Value = $handedBack
$wsZhCn.Range("K3").Value = "2016-10-17 15:58:58"

$wsZhCn.Range("C4").Value = $handedBack
$wsZhCn.Range("H4").Value = "2016-10-17 15:58:58"
$wsZhCn.Range("K4").Value = "2016-10-17 15:58:58"

# ---------------------------------------------------------------------
# de-de sheet: row 3 (7de06e11-...) and row 4 (a5786877-..., content
# duplicate of row 3) now have a handback target/file recorded.
# ---------------------------------------------------------------------
$wsDeDe.Range("C3").Value = $handedBack
$wsDeDe.Range("K3").Value = "2016-10-17 15:59:35"

$wsDeDe.Range("C4").Value = $handedBack
$wsDeDe.Range("K4").Value = "2016-10-17 15:59:35"

# ---------------------------------------------------------------------
# Re-create hyperlinks on both language sheets in final document order
# so that the new "Latest Target File" hyperlinks for I3/I4 are inserted
# alongside the pre-existing ones, and set the target/handback file text.
# ---------------------------------------------------------------------
function Rebuild-Hyperlinks($ws, $rId4Target, $rId3Target, $i3Display, $i3Target, $i4Display, $i4Target) {
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/946952792914813b00864f0d68375d0ffd7fae19/e2e/27c80a73-16fb-4437-a628-5ab6f9ace938.md", [Type]::Missing, [Type]::Missing, "27c80a73-16fb-4437-a628-5ab6f9ace938.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $rId3Target, [Type]::Missing, [Type]::Missing, "27c80a73-16fb-4437-a628-5ab6f9ace938.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3aed56102495107086df1b8044f434352684d4cb/e2e/7de06e11-3939-4d03-b6a1-b5617c4c3b14.md", [Type]::Missing, [Type]::Missing, "7de06e11-3939-4d03-b6a1-b5617c4c3b14.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $i3Target, [Type]::Missing, [Type]::Missing, $i3Display)
    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3aed56102495107086df1b8044f434352684d4cb/e2e/a5786877-59b2-44c3-a565-7c29da864af7.md", [Type]::Missing, [Type]::Missing, "a5786877-59b2-44c3-a565-7c29da864af7.md")
    $ws.Hyperlinks.Add($ws.Range("I4"), $i4Target, [Type]::Missing, [Type]::Missing, $i4Display)
}

Rebuild-Hyperlinks $wsZhCn `
    "" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fc2697803d62db337f8315c58d3867e51b2e4121/e2e/27c80a73-16fb-4437-a628-5ab6f9ace938.md" `
    "7de06e11-3939-4d03-b6a1-b5617c4c3b14.md" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3aed56102495107086df1b8044f434352684d4cb/e2e/7de06e11-3939-4d03-b6a1-b5617c4c3b14.md" `
    "7de06e11-3939-4d03-b6a1-b5617c4c3b14.md" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3aed56102495107086df1b8044f434352684d4cb/e2e/7de06e11-3939-4d03-b6a1-b5617c4c3b14.md"

Rebuild-Hyperlinks $wsDeDe `
    "" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0b147c7c4de9770f732f2d8e11be1046b2448ca1/e2e/27c80a73-16fb-4437-a628-5ab6f9ace938.md" `
    "7de06e11-3939-4d03-b6a1-b5617c4c3b14.md" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3aed56102495107086df1b8044f434352684d4cb/e2e/7de06e11-3939-4d03-b6a1-b5617c4c3b14.md" `
    "7de06e11-3939-4d03-b6a1-b5617c4c3b14.md" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3aed56102495107086df1b8044f434352684d4cb/e2e/7de06e11-3939-4d03-b6a1-b5617c4c3b14.md"

# Set the target/handback file names now that the hyperlinks (and their
# cell values) have been (re)written.
$wsZhCn.Range("J3").Value = "7de06e11-3939-4d03-b6a1-b5617c4c3b14.f6a6a3ae5ba24bd9a1e3ff558072adf0faa8e849.zh-cn.xlf"
$wsZhCn.Range("J4").Value = "7de06e11-3939-4d03-b6a1-b5617c4c3b14.f6a6a3ae5ba24bd9a1e3ff558072adf0faa8e849.zh-cn.xlf"

$wsDeDe.Range("J3").Value = "7de06e11-3939-4d03-b6a1-b5617c4c3b14.f6a6a3ae5ba24bd9a1e3ff558072adf0faa8e849.de-de.xlf"
$wsDeDe.Range("J4").Value = "7de06e11-3939-4d03-b6a1-b5617c4c3b14.f6a6a3ae5ba24bd9a1e3ff558072adf0faa8e849.de-de.xlf"
